$wb = $excel.ActiveWorkbook

# "Generate Report for Handoff": refresh the "Latest Handoff Datetime" for the
# 14432ec2-975e-438b-aa83-997f69c30a47 source file in both locale sheets.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-03-10 03:36:38"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-03-10 03:36:46"
